$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix isolated count increments in earlier rows ---
$ws.Cells.Item(229, 3).Value = 20
$ws.Cells.Item(239, 3).Value = 21
$ws.Cells.Item(1124, 3).Value = 18
$ws.Cells.Item(1190, 3).Value = 17
$ws.Cells.Item(1437, 3).Value = 7
$ws.Cells.Item(1512, 3).Value = 3
$ws.Cells.Item(1521, 3).Value = 5
$ws.Cells.Item(1539, 3).Value = 3
$ws.Cells.Item(1562, 3).Value = 2

# --- Rewrite date-shifted / expanded tail of the table (rows 1563-1662) ---
$ws.Cells.Item(1563, 1).Value = 44279
$ws.Cells.Item(1563, 2).Value = "50-59"
$ws.Cells.Item(1563, 3).Value = 1
$ws.Cells.Item(1564, 1).Value = 44279
$ws.Cells.Item(1564, 2).Value = "60-69"
$ws.Cells.Item(1564, 3).Value = 1
$ws.Cells.Item(1565, 1).Value = 44279
$ws.Cells.Item(1565, 2).Value = "70-79"
$ws.Cells.Item(1565, 3).Value = 3
$ws.Cells.Item(1566, 1).Value = 44279
$ws.Cells.Item(1566, 2).Value = "80+"
$ws.Cells.Item(1566, 3).Value = 2
$ws.Cells.Item(1567, 1).Value = 44280
$ws.Cells.Item(1567, 2).Value = "60-69"
$ws.Cells.Item(1567, 3).Value = 1
$ws.Cells.Item(1568, 1).Value = 44280
$ws.Cells.Item(1568, 2).Value = "70-79"
$ws.Cells.Item(1568, 3).Value = 6
$ws.Cells.Item(1569, 1).Value = 44280
$ws.Cells.Item(1569, 2).Value = "80+"
$ws.Cells.Item(1569, 3).Value = 2
$ws.Cells.Item(1570, 1).Value = 44281
$ws.Cells.Item(1570, 2).Value = "30-39"
$ws.Cells.Item(1570, 3).Value = 1
$ws.Cells.Item(1571, 1).Value = 44281
$ws.Cells.Item(1571, 2).Value = "50-59"
$ws.Cells.Item(1571, 3).Value = 1
$ws.Cells.Item(1572, 1).Value = 44281
$ws.Cells.Item(1572, 2).Value = "60-69"
$ws.Cells.Item(1572, 3).Value = 3
$ws.Cells.Item(1573, 1).Value = 44281
$ws.Cells.Item(1573, 2).Value = "70-79"
$ws.Cells.Item(1573, 3).Value = 1
$ws.Cells.Item(1574, 1).Value = 44281
$ws.Cells.Item(1574, 2).Value = "80+"
$ws.Cells.Item(1574, 3).Value = 2
$ws.Cells.Item(1575, 1).Value = 44282
$ws.Cells.Item(1575, 2).Value = "60-69"
$ws.Cells.Item(1575, 3).Value = 1
$ws.Cells.Item(1576, 1).Value = 44282
$ws.Cells.Item(1576, 2).Value = "70-79"
$ws.Cells.Item(1576, 3).Value = 1
$ws.Cells.Item(1577, 1).Value = 44282
$ws.Cells.Item(1577, 2).Value = "80+"
$ws.Cells.Item(1577, 3).Value = 5
$ws.Cells.Item(1578, 1).Value = 44283
$ws.Cells.Item(1578, 2).Value = "40-49"
$ws.Cells.Item(1578, 3).Value = 1
$ws.Cells.Item(1579, 1).Value = 44283
$ws.Cells.Item(1579, 2).Value = "50-59"
$ws.Cells.Item(1579, 3).Value = 1
$ws.Cells.Item(1580, 1).Value = 44283
$ws.Cells.Item(1580, 2).Value = "60-69"
$ws.Cells.Item(1580, 3).Value = 1
$ws.Cells.Item(1581, 1).Value = 44283
$ws.Cells.Item(1581, 2).Value = "70-79"
$ws.Cells.Item(1581, 3).Value = 4
$ws.Cells.Item(1582, 1).Value = 44283
$ws.Cells.Item(1582, 2).Value = "80+"
$ws.Cells.Item(1582, 3).Value = 1
$ws.Cells.Item(1583, 1).Value = 44284
$ws.Cells.Item(1583, 2).Value = "50-59"
$ws.Cells.Item(1583, 3).Value = 1
$ws.Cells.Item(1584, 1).Value = 44284
$ws.Cells.Item(1584, 2).Value = "60-69"
$ws.Cells.Item(1584, 3).Value = 1
$ws.Cells.Item(1585, 1).Value = 44284
$ws.Cells.Item(1585, 2).Value = "70-79"
$ws.Cells.Item(1585, 3).Value = 2
$ws.Cells.Item(1586, 1).Value = 44284
$ws.Cells.Item(1586, 2).Value = "80+"
$ws.Cells.Item(1586, 3).Value = 1
$ws.Cells.Item(1587, 1).Value = 44285
$ws.Cells.Item(1587, 2).Value = "60-69"
$ws.Cells.Item(1587, 3).Value = 5
$ws.Cells.Item(1588, 1).Value = 44285
$ws.Cells.Item(1588, 2).Value = "80+"
$ws.Cells.Item(1588, 3).Value = 3
$ws.Cells.Item(1589, 1).Value = 44286
$ws.Cells.Item(1589, 2).Value = "40-49"
$ws.Cells.Item(1589, 3).Value = 1
$ws.Cells.Item(1590, 1).Value = 44286
$ws.Cells.Item(1590, 2).Value = "50-59"
$ws.Cells.Item(1590, 3).Value = 3
$ws.Cells.Item(1591, 1).Value = 44286
$ws.Cells.Item(1591, 2).Value = "60-69"
$ws.Cells.Item(1591, 3).Value = 1
$ws.Cells.Item(1592, 1).Value = 44286
$ws.Cells.Item(1592, 2).Value = "70-79"
$ws.Cells.Item(1592, 3).Value = 2
$ws.Cells.Item(1593, 1).Value = 44286
$ws.Cells.Item(1593, 2).Value = "80+"
$ws.Cells.Item(1593, 3).Value = 2
$ws.Cells.Item(1594, 1).Value = 44287
$ws.Cells.Item(1594, 2).Value = "50-59"
$ws.Cells.Item(1594, 3).Value = 2
$ws.Cells.Item(1595, 1).Value = 44287
$ws.Cells.Item(1595, 2).Value = "60-69"
$ws.Cells.Item(1595, 3).Value = 1
$ws.Cells.Item(1596, 1).Value = 44287
$ws.Cells.Item(1596, 2).Value = "70-79"
$ws.Cells.Item(1596, 3).Value = 1
$ws.Cells.Item(1597, 1).Value = 44287
$ws.Cells.Item(1597, 2).Value = "80+"
$ws.Cells.Item(1597, 3).Value = 3
$ws.Cells.Item(1598, 1).Value = 44288
$ws.Cells.Item(1598, 2).Value = "40-49"
$ws.Cells.Item(1598, 3).Value = 1
$ws.Cells.Item(1599, 1).Value = 44288
$ws.Cells.Item(1599, 2).Value = "60-69"
$ws.Cells.Item(1599, 3).Value = 1
$ws.Cells.Item(1600, 1).Value = 44288
$ws.Cells.Item(1600, 2).Value = "70-79"
$ws.Cells.Item(1600, 3).Value = 3
$ws.Cells.Item(1601, 1).Value = 44288
$ws.Cells.Item(1601, 2).Value = "80+"
$ws.Cells.Item(1601, 3).Value = 1
$ws.Cells.Item(1602, 1).Value = 44289
$ws.Cells.Item(1602, 2).Value = "60-69"
$ws.Cells.Item(1602, 3).Value = 1
$ws.Cells.Item(1603, 1).Value = 44289
$ws.Cells.Item(1603, 2).Value = "70-79"
$ws.Cells.Item(1603, 3).Value = 4
$ws.Cells.Item(1604, 1).Value = 44289
$ws.Cells.Item(1604, 2).Value = "80+"
$ws.Cells.Item(1604, 3).Value = 1
$ws.Cells.Item(1605, 1).Value = 44290
$ws.Cells.Item(1605, 2).Value = "30-39"
$ws.Cells.Item(1605, 3).Value = 1
$ws.Cells.Item(1606, 1).Value = 44290
$ws.Cells.Item(1606, 2).Value = "70-79"
$ws.Cells.Item(1606, 3).Value = 2
$ws.Cells.Item(1607, 1).Value = 44290
$ws.Cells.Item(1607, 2).Value = "80+"
$ws.Cells.Item(1607, 3).Value = 2
$ws.Cells.Item(1608, 1).Value = 44291
$ws.Cells.Item(1608, 2).Value = "80+"
$ws.Cells.Item(1608, 3).Value = 2
$ws.Cells.Item(1609, 1).Value = 44292
$ws.Cells.Item(1609, 2).Value = "50-59"
$ws.Cells.Item(1609, 3).Value = 2
$ws.Cells.Item(1610, 1).Value = 44292
$ws.Cells.Item(1610, 2).Value = "60-69"
$ws.Cells.Item(1610, 3).Value = 2
$ws.Cells.Item(1611, 1).Value = 44292
$ws.Cells.Item(1611, 2).Value = "70-79"
$ws.Cells.Item(1611, 3).Value = 1
$ws.Cells.Item(1612, 1).Value = 44293
$ws.Cells.Item(1612, 2).Value = "60-69"
$ws.Cells.Item(1612, 3).Value = 2
$ws.Cells.Item(1613, 1).Value = 44293
$ws.Cells.Item(1613, 2).Value = "70-79"
$ws.Cells.Item(1613, 3).Value = 1
$ws.Cells.Item(1614, 1).Value = 44293
$ws.Cells.Item(1614, 2).Value = "80+"
$ws.Cells.Item(1614, 3).Value = 4
$ws.Cells.Item(1615, 1).Value = 44294
$ws.Cells.Item(1615, 2).Value = "50-59"
$ws.Cells.Item(1615, 3).Value = 1
$ws.Cells.Item(1616, 1).Value = 44294
$ws.Cells.Item(1616, 2).Value = "70-79"
$ws.Cells.Item(1616, 3).Value = 3
$ws.Cells.Item(1617, 1).Value = 44294
$ws.Cells.Item(1617, 2).Value = "80+"
$ws.Cells.Item(1617, 3).Value = 2
$ws.Cells.Item(1618, 1).Value = 44295
$ws.Cells.Item(1618, 2).Value = "50-59"
$ws.Cells.Item(1618, 3).Value = 1
$ws.Cells.Item(1619, 1).Value = 44295
$ws.Cells.Item(1619, 2).Value = "60-69"
$ws.Cells.Item(1619, 3).Value = 3
$ws.Cells.Item(1620, 1).Value = 44295
$ws.Cells.Item(1620, 2).Value = "70-79"
$ws.Cells.Item(1620, 3).Value = 2
$ws.Cells.Item(1621, 1).Value = 44295
$ws.Cells.Item(1621, 2).Value = "80+"
$ws.Cells.Item(1621, 3).Value = 1
$ws.Cells.Item(1622, 1).Value = 44296
$ws.Cells.Item(1622, 2).Value = "60-69"
$ws.Cells.Item(1622, 3).Value = 1
$ws.Cells.Item(1623, 1).Value = 44296
$ws.Cells.Item(1623, 2).Value = "70-79"
$ws.Cells.Item(1623, 3).Value = 3
$ws.Cells.Item(1624, 1).Value = 44296
$ws.Cells.Item(1624, 2).Value = "80+"
$ws.Cells.Item(1624, 3).Value = 4
$ws.Cells.Item(1625, 1).Value = 44297
$ws.Cells.Item(1625, 2).Value = "40-49"
$ws.Cells.Item(1625, 3).Value = 1
$ws.Cells.Item(1626, 1).Value = 44297
$ws.Cells.Item(1626, 2).Value = "60-69"
$ws.Cells.Item(1626, 3).Value = 2
$ws.Cells.Item(1627, 1).Value = 44297
$ws.Cells.Item(1627, 2).Value = "70-79"
$ws.Cells.Item(1627, 3).Value = 1
$ws.Cells.Item(1628, 1).Value = 44297
$ws.Cells.Item(1628, 2).Value = "80+"
$ws.Cells.Item(1628, 3).Value = 4
$ws.Cells.Item(1629, 1).Value = 44298
$ws.Cells.Item(1629, 2).Value = "40-49"
$ws.Cells.Item(1629, 3).Value = 1
$ws.Cells.Item(1630, 1).Value = 44298
$ws.Cells.Item(1630, 2).Value = "50-59"
$ws.Cells.Item(1630, 3).Value = 2
$ws.Cells.Item(1631, 1).Value = 44298
$ws.Cells.Item(1631, 2).Value = "60-69"
$ws.Cells.Item(1631, 3).Value = 4
$ws.Cells.Item(1632, 1).Value = 44298
$ws.Cells.Item(1632, 2).Value = "70-79"
$ws.Cells.Item(1632, 3).Value = 2
$ws.Cells.Item(1633, 1).Value = 44298
$ws.Cells.Item(1633, 2).Value = "80+"
$ws.Cells.Item(1633, 3).Value = 1
$ws.Cells.Item(1634, 1).Value = 44299
$ws.Cells.Item(1634, 2).Value = "40-49"
$ws.Cells.Item(1634, 3).Value = 1
$ws.Cells.Item(1635, 1).Value = 44299
$ws.Cells.Item(1635, 2).Value = "50-59"
$ws.Cells.Item(1635, 3).Value = 1
$ws.Cells.Item(1636, 1).Value = 44299
$ws.Cells.Item(1636, 2).Value = "60-69"
$ws.Cells.Item(1636, 3).Value = 5
$ws.Cells.Item(1637, 1).Value = 44299
$ws.Cells.Item(1637, 2).Value = "80+"
$ws.Cells.Item(1637, 3).Value = 6
$ws.Cells.Item(1638, 1).Value = 44300
$ws.Cells.Item(1638, 2).Value = "50-59"
$ws.Cells.Item(1638, 3).Value = 1
$ws.Cells.Item(1639, 1).Value = 44300
$ws.Cells.Item(1639, 2).Value = "60-69"
$ws.Cells.Item(1639, 3).Value = 1
$ws.Cells.Item(1640, 1).Value = 44300
$ws.Cells.Item(1640, 2).Value = "70-79"
$ws.Cells.Item(1640, 3).Value = 1
$ws.Cells.Item(1641, 1).Value = 44301
$ws.Cells.Item(1641, 2).Value = "60-69"
$ws.Cells.Item(1641, 3).Value = 1
$ws.Cells.Item(1642, 1).Value = 44301
$ws.Cells.Item(1642, 2).Value = "70-79"
$ws.Cells.Item(1642, 3).Value = 1
$ws.Cells.Item(1643, 1).Value = 44301
$ws.Cells.Item(1643, 2).Value = "80+"
$ws.Cells.Item(1643, 3).Value = 2
$ws.Cells.Item(1644, 1).Value = 44302
$ws.Cells.Item(1644, 2).Value = "30-39"
$ws.Cells.Item(1644, 3).Value = 2
$ws.Cells.Item(1645, 1).Value = 44302
$ws.Cells.Item(1645, 2).Value = "40-49"
$ws.Cells.Item(1645, 3).Value = 1
$ws.Cells.Item(1646, 1).Value = 44302
$ws.Cells.Item(1646, 2).Value = "60-69"
$ws.Cells.Item(1646, 3).Value = 1
$ws.Cells.Item(1647, 1).Value = 44302
$ws.Cells.Item(1647, 2).Value = "70-79"
$ws.Cells.Item(1647, 3).Value = 4
$ws.Cells.Item(1648, 1).Value = 44302
$ws.Cells.Item(1648, 2).Value = "80+"
$ws.Cells.Item(1648, 3).Value = 1
$ws.Cells.Item(1649, 1).Value = 44303
$ws.Cells.Item(1649, 2).Value = "60-69"
$ws.Cells.Item(1649, 3).Value = 1
$ws.Cells.Item(1650, 1).Value = 44303
$ws.Cells.Item(1650, 2).Value = "70-79"
$ws.Cells.Item(1650, 3).Value = 1
$ws.Cells.Item(1651, 1).Value = 44303
$ws.Cells.Item(1651, 2).Value = "80+"
$ws.Cells.Item(1651, 3).Value = 1
$ws.Cells.Item(1652, 1).Value = 44304
$ws.Cells.Item(1652, 2).Value = "60-69"
$ws.Cells.Item(1652, 3).Value = 1
$ws.Cells.Item(1653, 1).Value = 44304
$ws.Cells.Item(1653, 2).Value = "80+"
$ws.Cells.Item(1653, 3).Value = 2
$ws.Cells.Item(1654, 1).Value = 44305
$ws.Cells.Item(1654, 2).Value = "60-69"
$ws.Cells.Item(1654, 3).Value = 1
$ws.Cells.Item(1655, 1).Value = 44305
$ws.Cells.Item(1655, 2).Value = "70-79"
$ws.Cells.Item(1655, 3).Value = 1
$ws.Cells.Item(1656, 1).Value = 44305
$ws.Cells.Item(1656, 2).Value = "80+"
$ws.Cells.Item(1656, 3).Value = 1
$ws.Cells.Item(1657, 1).Value = 44306
$ws.Cells.Item(1657, 2).Value = "40-49"
$ws.Cells.Item(1657, 3).Value = 1
$ws.Cells.Item(1658, 1).Value = 44306
$ws.Cells.Item(1658, 2).Value = "60-69"
$ws.Cells.Item(1658, 3).Value = 2
$ws.Cells.Item(1659, 1).Value = 44306
$ws.Cells.Item(1659, 2).Value = "80+"
$ws.Cells.Item(1659, 3).Value = 2
$ws.Cells.Item(1660, 1).Value = 44307
$ws.Cells.Item(1660, 2).Value = "80+"
$ws.Cells.Item(1660, 3).Value = 2
$ws.Cells.Item(1661, 1).Value = 44308
$ws.Cells.Item(1661, 2).Value = "50-59"
$ws.Cells.Item(1661, 3).Value = 1
$ws.Cells.Item(1662, 1).Value = 44308
$ws.Cells.Item(1662, 2).Value = "80+"
$ws.Cells.Item(1662, 3).Value = 1

# Ensure new date cells (rows beyond the original 1631) use the date number format
$ws.Range("A1632:A1662").NumberFormat = "YYYY-MM-DD HH:MM:SS"
